$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 149 (CHO-Medimix Sandal with Eladi oils for g)
$ws.Cells.Item(149, 2).Value = 63902
$ws.Cells.Item(149, 3).Value = 'CHO-Medimix Sandal with Eladi oils for glowing skin and natural protection Soap-125gms'
$ws.Cells.Item(149, 4).Value = 32.02
$ws.Cells.Item(149, 5).Value = 34.04
$ws.Cells.Item(149, 6).Value = 2
$ws.Cells.Item(149, 7).Value = 64.04000000000001

# Row 150 (CHO-Medimix Sandal with Eladi oils for g)
$ws.Cells.Item(150, 2).Value = 48654
$ws.Cells.Item(150, 3).Value = 'CHO-Medimix Sandal with Eladi oils for glowing skin and natural protection Soap-125gms'
$ws.Cells.Item(150, 4).Value = 32.02
$ws.Cells.Item(150, 5).Value = 38.26
$ws.Cells.Item(150, 6).Value = -1
$ws.Cells.Item(150, 7).Value = -32.02

# Row 161 (COL-Colgate Zigzag Charcoal Pack of 4 To)
$ws.Cells.Item(161, 2).Value = 53925
$ws.Cells.Item(161, 3).Value = 'COL-Colgate Zigzag Charcoal Pack of 4 Toothbrush'
$ws.Cells.Item(161, 4).Value = 66.44
$ws.Cells.Item(161, 5).Value = 79.37
$ws.Cells.Item(161, 6).Value = 1
$ws.Cells.Item(161, 7).Value = 66.44

# Row 162 (COL-Colgate Zigzag Charcoal Pack of 4 To)
$ws.Cells.Item(162, 2).Value = 64350
$ws.Cells.Item(162, 3).Value = 'COL-Colgate Zigzag Charcoal Pack of 4 Toothbrush'
$ws.Cells.Item(162, 4).Value = 66.44
$ws.Cells.Item(162, 5).Value = 70.63
$ws.Cells.Item(162, 6).Value = 29
$ws.Cells.Item(162, 7).Value = 1926.76

# Row 163 (COL-Colgate Zigzag Charcoal Pack of 4 To)
$ws.Cells.Item(163, 2).Value = 57756
$ws.Cells.Item(163, 3).Value = 'COL-Colgate Zigzag Charcoal Pack of 4 Toothbrush'
$ws.Cells.Item(163, 4).Value = 66.44
$ws.Cells.Item(163, 5).Value = 79.37
$ws.Cells.Item(163, 6).Value = -100
$ws.Cells.Item(163, 7).Value = -6644

# Row 316 (HUL-Bru Inst Poly 50g)
$ws.Cells.Item(316, 2).Value = 61610
$ws.Cells.Item(316, 3).Value = 'HUL-Bru Inst Poly 50g'
$ws.Cells.Item(316, 4).Value = 102.71
$ws.Cells.Item(316, 5).Value = 122.71
$ws.Cells.Item(316, 6).Value = -58
$ws.Cells.Item(316, 7).Value = -5957.18

# Row 317 (HUL-Bru Inst Poly 50g)
$ws.Cells.Item(317, 2).Value = 57077
$ws.Cells.Item(317, 3).Value = 'HUL-Bru Inst Poly 50g'
$ws.Cells.Item(317, 4).Value = 93.08
$ws.Cells.Item(317, 5).Value = 111.2
$ws.Cells.Item(317, 6).Value = 1
$ws.Cells.Item(317, 7).Value = 93.08

# Row 318 (HUL-Bru Inst Poly 50g)
$ws.Cells.Item(318, 2).Value = 63565
$ws.Cells.Item(318, 3).Value = 'HUL-Bru Inst Poly 50g'
$ws.Cells.Item(318, 4).Value = 102.71
$ws.Cells.Item(318, 5).Value = 109.19
$ws.Cells.Item(318, 6).Value = 60
$ws.Cells.Item(318, 7).Value = 6162.6

# Row 346 (HUL-Kissan nango jam 490g)
$ws.Cells.Item(346, 2).Value = 55373
$ws.Cells.Item(346, 3).Value = 'HUL-Kissan nango jam 490g'
$ws.Cells.Item(346, 4).Value = 144.28
$ws.Cells.Item(346, 5).Value = 163.62
$ws.Cells.Item(346, 6).Value = -94
$ws.Cells.Item(346, 7).Value = -13562.32

# Row 347 (HUL-Kissan nango jam 490g)
$ws.Cells.Item(347, 2).Value = 63520
$ws.Cells.Item(347, 3).Value = 'HUL-Kissan nango jam 490g'
$ws.Cells.Item(347, 4).Value = 144.28
$ws.Cells.Item(347, 5).Value = 153.4
$ws.Cells.Item(347, 6).Value = 91
$ws.Cells.Item(347, 7).Value = 13129.48

# Row 350 (HUL-Kissan Pineapple Jam 500G)
$ws.Cells.Item(350, 2).Value = 57802
$ws.Cells.Item(350, 3).Value = 'HUL-Kissan Pineapple Jam 500G'
$ws.Cells.Item(350, 4).Value = 143.48
$ws.Cells.Item(350, 5).Value = 162.71
$ws.Cells.Item(350, 6).Value = -79
$ws.Cells.Item(350, 7).Value = -11334.92

# Row 351 (HUL-Kissan Pineapple Jam 500G)
$ws.Cells.Item(351, 2).Value = 63571
$ws.Cells.Item(351, 3).Value = 'HUL-Kissan Pineapple Jam 500G'
$ws.Cells.Item(351, 4).Value = 143.48
$ws.Cells.Item(351, 5).Value = 152.53
$ws.Cells.Item(351, 6).Value = 18
$ws.Cells.Item(351, 7).Value = 2582.64

# Row 352 (HUL-Kissan Pineapple Jam 500G)
$ws.Cells.Item(352, 2).Value = 63531
$ws.Cells.Item(352, 3).Value = 'HUL-Kissan Pineapple Jam 500G'
$ws.Cells.Item(352, 4).Value = 143.48
$ws.Cells.Item(352, 5).Value = 152.53
$ws.Cells.Item(352, 6).Value = 80
$ws.Cells.Item(352, 7).Value = 11478.4

# Row 372 (HUL-Liril Soap 125 G)
$ws.Cells.Item(372, 2).Value = 63652
$ws.Cells.Item(372, 3).Value = 'HUL-Liril Soap 125 G'
$ws.Cells.Item(372, 4).Value = 52.13
$ws.Cells.Item(372, 5).Value = 55.42
$ws.Cells.Item(372, 6).Value = 192
$ws.Cells.Item(372, 7).Value = 10008.96

# Row 373 (HUL-Liril Soap 125 G)
$ws.Cells.Item(373, 2).Value = 57885
$ws.Cells.Item(373, 3).Value = 'HUL-Liril Soap 125 G'
$ws.Cells.Item(373, 4).Value = 52.13
$ws.Cells.Item(373, 5).Value = 62.28
$ws.Cells.Item(373, 6).Value = 4
$ws.Cells.Item(373, 7).Value = 208.52

# Row 375 (HUL-lux advanced eventoned glow 4x100)
$ws.Cells.Item(375, 2).Value = 63563
$ws.Cells.Item(375, 3).Value = 'HUL-lux advanced eventoned glow 4x100'
$ws.Cells.Item(375, 4).Value = 111.96
$ws.Cells.Item(375, 5).Value = 119.04
$ws.Cells.Item(375, 6).Value = 2
$ws.Cells.Item(375, 7).Value = 223.92

# Row 376 (HUL-lux advanced eventoned glow 4x100)
$ws.Cells.Item(376, 2).Value = 61605
$ws.Cells.Item(376, 3).Value = 'HUL-lux advanced eventoned glow 4x100'
$ws.Cells.Item(376, 4).Value = 111.96
$ws.Cells.Item(376, 5).Value = 133.78
$ws.Cells.Item(376, 6).Value = -13
$ws.Cells.Item(376, 7).Value = -1455.48

# Row 382 (Hul-pears pure and gentle 3x125 gm)
$ws.Cells.Item(382, 2).Value = 63560
$ws.Cells.Item(382, 3).Value = 'Hul-pears pure and gentle 3x125 gm'
$ws.Cells.Item(382, 4).Value = 126.86
$ws.Cells.Item(382, 5).Value = 134.87
$ws.Cells.Item(382, 6).Value = 14
$ws.Cells.Item(382, 7).Value = 1776.04

# Row 383 (Hul-pears pure and gentle 3x125 gm)
$ws.Cells.Item(383, 2).Value = 60325
$ws.Cells.Item(383, 3).Value = 'Hul-pears pure and gentle 3x125 gm'
$ws.Cells.Item(383, 4).Value = 126.86
$ws.Cells.Item(383, 5).Value = 151.57
$ws.Cells.Item(383, 6).Value = -102
$ws.Cells.Item(383, 7).Value = -12939.72

# Row 389 (HUL-Rap Refresh Bolt 1Kg)
$ws.Cells.Item(389, 2).Value = 57817
$ws.Cells.Item(389, 3).Value = 'HUL-Rap Refresh Bolt 1Kg'
$ws.Cells.Item(389, 4).Value = 79.81
$ws.Cells.Item(389, 5).Value = 95.34999999999999
$ws.Cells.Item(389, 6).Value = 3
$ws.Cells.Item(389, 7).Value = 239.43

# Row 390 (HUL-Rap Refresh Bolt 1Kg)
$ws.Cells.Item(390, 2).Value = 62865
$ws.Cells.Item(390, 3).Value = 'HUL-Rap Refresh Bolt 1Kg'
$ws.Cells.Item(390, 4).Value = 79.81
$ws.Cells.Item(390, 5).Value = 95.34999999999999
$ws.Cells.Item(390, 6).Value = 33
$ws.Cells.Item(390, 7).Value = 2633.73

# Row 400 (HUL-Sfxl Ew Bale 500G)
$ws.Cells.Item(400, 2).Value = 62933
$ws.Cells.Item(400, 3).Value = 'HUL-Sfxl Ew Bale 500G'
$ws.Cells.Item(400, 4).Value = 59.13
$ws.Cells.Item(400, 5).Value = 70.65000000000001
$ws.Cells.Item(400, 6).Value = 129
$ws.Cells.Item(400, 7).Value = 7627.77

# Row 401 (HUL-Sfxl Ew Bale 500G)
$ws.Cells.Item(401, 2).Value = 57835
$ws.Cells.Item(401, 3).Value = 'HUL-Sfxl Ew Bale 500G'
$ws.Cells.Item(401, 4).Value = 59.13
$ws.Cells.Item(401, 5).Value = 70.65000000000001
$ws.Cells.Item(401, 6).Value = 1
$ws.Cells.Item(401, 7).Value = 59.13

# Row 419 (HUL-Surf Exl Mtc Liq Fl 1 Ltr Cp)
$ws.Cells.Item(419, 2).Value = 57856
$ws.Cells.Item(419, 3).Value = 'HUL-Surf Exl Mtc Liq Fl 1 Ltr Cp'
$ws.Cells.Item(419, 4).Value = 171.33
$ws.Cells.Item(419, 5).Value = 204.69
$ws.Cells.Item(419, 6).Value = 2
$ws.Cells.Item(419, 7).Value = 342.66

# Row 420 (HUL-Surf Exl Mtc Liq Fl 1 Ltr Cp)
$ws.Cells.Item(420, 2).Value = 63007
$ws.Cells.Item(420, 3).Value = 'HUL-Surf Exl Mtc Liq Fl 1 Ltr Cp'
$ws.Cells.Item(420, 4).Value = 171.33
$ws.Cells.Item(420, 5).Value = 204.69
$ws.Cells.Item(420, 6).Value = 852
$ws.Cells.Item(420, 7).Value = 145973.16

# Row 431 (HUL-Vim Bar Multipack Fw 4X200G)
$ws.Cells.Item(431, 2).Value = 63102
$ws.Cells.Item(431, 3).Value = 'HUL-Vim Bar Multipack Fw 4X200G'
$ws.Cells.Item(431, 4).Value = 59.47
$ws.Cells.Item(431, 5).Value = 71.05
$ws.Cells.Item(431, 6).Value = 4
$ws.Cells.Item(431, 7).Value = 237.88

# Row 432 (HUL-VIM BAR MULTIPACK FW 4X200G)
$ws.Cells.Item(432, 2).Value = 53082
$ws.Cells.Item(432, 3).Value = 'HUL-VIM BAR MULTIPACK FW 4X200G'
$ws.Cells.Item(432, 4).Value = 59.47
$ws.Cells.Item(432, 5).Value = 71.05
$ws.Cells.Item(432, 6).Value = 1
$ws.Cells.Item(432, 7).Value = 59.47

# Row 586 (CRE-Cremica Chocolate Cream 150Gm)
$ws.Cells.Item(586, 2).Value = 45695
$ws.Cells.Item(586, 3).Value = 'CRE-Cremica Chocolate Cream 150Gm'
$ws.Cells.Item(586, 4).Value = 19.73
$ws.Cells.Item(586, 5).Value = 23.58
$ws.Cells.Item(586, 6).Value = -36
$ws.Cells.Item(586, 7).Value = -710.28

# Row 587 (CRE-Cremica Chocolate Cream 150Gm)
$ws.Cells.Item(587, 2).Value = 64915
$ws.Cells.Item(587, 3).Value = 'CRE-Cremica Chocolate Cream 150Gm'
$ws.Cells.Item(587, 4).Value = 19.73
$ws.Cells.Item(587, 5).Value = 20.98
$ws.Cells.Item(587, 6).Value = 13
$ws.Cells.Item(587, 7).Value = 256.49

# Row 599 (CRE-Cremica Oatmeal Digestive 112.5 Gm)
$ws.Cells.Item(599, 2).Value = 64925
$ws.Cells.Item(599, 3).Value = 'CRE-Cremica Oatmeal Digestive 112.5 Gm'
$ws.Cells.Item(599, 4).Value = 13.15
$ws.Cells.Item(599, 5).Value = 13.97
$ws.Cells.Item(599, 6).Value = 273
$ws.Cells.Item(599, 7).Value = 3589.95

# Row 600 (CRE-Cremica Oatmeal Digestive 112.5 Gm)
$ws.Cells.Item(600, 2).Value = 45709
$ws.Cells.Item(600, 3).Value = 'CRE-Cremica Oatmeal Digestive 112.5 Gm'
$ws.Cells.Item(600, 4).Value = 13.15
$ws.Cells.Item(600, 5).Value = 15.69
$ws.Cells.Item(600, 6).Value = -300
$ws.Cells.Item(600, 7).Value = -3945

# Row 687 (PRI-B-50 VIMAL Copper Glass 300ML (2pc S)
$ws.Cells.Item(687, 2).Value = 64810
$ws.Cells.Item(687, 3).Value = 'PRI-B-50 VIMAL Copper Glass 300ML (2pc Set)'
$ws.Cells.Item(687, 4).Value = 273.92
$ws.Cells.Item(687, 5).Value = 291.22
$ws.Cells.Item(687, 6).Value = 7
$ws.Cells.Item(687, 7).Value = 1917.44

# Row 688 (PRI-B-50 VIMAL Copper Glass 300ML (2pc S)
$ws.Cells.Item(688, 2).Value = 53319
$ws.Cells.Item(688, 3).Value = 'PRI-B-50 VIMAL Copper Glass 300ML (2pc Set)'
$ws.Cells.Item(688, 4).Value = 273.92
$ws.Cells.Item(688, 5).Value = 310.64
$ws.Cells.Item(688, 6).Value = -6
$ws.Cells.Item(688, 7).Value = -1643.52

# Row 720 (Rasna Nagpur Orange (32 Glass))
$ws.Cells.Item(720, 2).Value = 64830
$ws.Cells.Item(720, 3).Value = 'Rasna Nagpur Orange (32 Glass)'
$ws.Cells.Item(720, 4).Value = 32.83
$ws.Cells.Item(720, 5).Value = 34.9
$ws.Cells.Item(720, 6).Value = 115
$ws.Cells.Item(720, 7).Value = 3775.45

# Row 721 (Rasna Nagpur Orange (32 Glass))
$ws.Cells.Item(721, 2).Value = 60022
$ws.Cells.Item(721, 3).Value = 'Rasna Nagpur Orange (32 Glass)'
$ws.Cells.Item(721, 4).Value = 32.83
$ws.Cells.Item(721, 5).Value = 37.22
$ws.Cells.Item(721, 6).Value = -113
$ws.Cells.Item(721, 7).Value = -3709.79
